$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

# Column A holds dates stored as plain text (shared strings) in this sheet,
# not real Excel dates. Force the cell to Text format before assigning the
# value so Excel doesn't auto-convert "2020-07-26" into a date serial, then
# restore the default "Normal" style so no stray style index is left on the
# cell (matching the rest of the data rows which carry no explicit style).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2020-07-26"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 390516
$ws.Cells.Item($row, 3).Value = 438468
$ws.Cells.Item($row, 4).Value = 89397
$ws.Cells.Item($row, 5).Value = 43680
$ws.Cells.Item($row, 6).Value = 27.78
